# Update ticket/attendance counts ("F" column) on all sheets to match
# the refreshed data snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 4193
$ws1.Range("F10").Value = 398
$ws1.Range("F11").Value = 3302
$ws1.Range("F12").Value = 918
$ws1.Range("F18").Value = 17
$ws1.Range("F20").Value = 459
$ws1.Range("F22").Value = 42
$ws1.Range("F23").Value = 9233
$ws1.Range("F24").Value = 5855
$ws1.Range("F26").Value = 195
$ws1.Range("F29").Value = 801
$ws1.Range("F34").Value = 84
$ws1.Range("F38").Value = 11
$ws1.Range("F39").Value = 942

# 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 71

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 8540
$ws3.Range("F3").Value = 396
$ws3.Range("F4").Value = 1437

# 全部类型 (All Types) - combined view of the three sheets above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 8540
$ws4.Range("F4").Value  = 396
$ws4.Range("F5").Value  = 1437
$ws4.Range("F10").Value = 4193
$ws4.Range("F12").Value = 398
$ws4.Range("F13").Value = 3302
$ws4.Range("F14").Value = 918
$ws4.Range("F25").Value = 17
$ws4.Range("F27").Value = 459
$ws4.Range("F29").Value = 42
$ws4.Range("F30").Value = 9233
$ws4.Range("F32").Value = 71
$ws4.Range("F34").Value = 195
$ws4.Range("F35").Value = 801
$ws4.Range("F38").Value = 84
$ws4.Range("F43").Value = 11
$ws4.Range("F44").Value = 942
